$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "halosalsa4" -> "halosalsa5" / "halosalsa4@gmail.com" -> "halosalsa5@gmail.com"
$ws.Range("B7").Value = "halosalsa5"
$ws.Range("B8").Value = "halosalsa5@gmail.com"

# Move the active cell/selection from B8 to H8
$ws.Activate()
$ws.Range("H8").Select()
